$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff appends a brand-new row 98 at the bottom of the (single) sheet,
# shaped just like the preceding row 97: a date in column A (using the
# existing date-formatted style), plain numeric values in columns B-F, and
# text columns G/H whose contents are already present elsewhere in the
# sheet (and therefore reuse existing shared-string entries instead of
# creating new ones).

$xlPasteAll = -4104
$xlPasteFormats = -4122

# Columns A-F: copy the existing number formatting from row 97 down to the
# new row 98, then fill in the new values.
$ws.Range("A97:F97").Copy()
$ws.Range("A98:F98").PasteSpecial($xlPasteFormats)

# Columns G (close price, stored as a numeric-looking shared string) and H
# (ticker) hold exactly the same content as row 97, so copy the whole
# cells (value + format) rather than retyping them - this keeps reusing
# the same shared-string entries / default style instead of creating new
# ones.
$ws.Range("G97:H97").Copy()
$ws.Range("G98:H98").PasteSpecial($xlPasteAll)

$ws.Cells.Item(98, 1).Value = 45463.2916666667
$ws.Cells.Item(98, 2).Value = 0
$ws.Cells.Item(98, 3).Value = 0.699999988079071
$ws.Cells.Item(98, 4).Value = 0.699999988079071
$ws.Cells.Item(98, 5).Value = 0.699999988079071
$ws.Cells.Item(98, 6).Value = 0.699999988079071

$excel.CutCopyMode = 0
